$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.19725239276886
$ws.Range("B1").Value = 2.347532749176025
$ws.Range("C1").Value = 6.828729629516602
$ws.Range("D1").Value = 2.322959423065186
$ws.Range("E1").Value = 1.183044075965881
